$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the cells that changed in existing rows (3,4,5,8,9,10) ---
$ws.Range("G3").Value = -20
$ws.Range("I3").Value = 0.26

$ws.Range("G4").Value = -250

$ws.Range("G5").Value = -109

$ws.Range("G8").Value = -26
$ws.Range("H8").Value = 1.07

$ws.Range("G9").Value = -57
$ws.Range("I9").Value = 0.17

$ws.Range("G10").Value = -130
$ws.Range("I10").Value = 1.46

# --- Append the new row 11 ---
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "2025-08-04"
$ws.Range("A11").ClearFormats()

$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "BEMOL S/A"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "391921"
$ws.Range("D11").ClearFormats()

$ws.Range("E11").Value = 10130
$ws.Range("F11").Value = "FONE DE OUVIDO SEM FIO A GOLD V5.3"
$ws.Range("G11").Value = -1248
$ws.Range("H11").Value = 1.06
$ws.Range("I11").Value = 0.31
